# Fruta / hortaliza, semanal
# A new daily price record (Naranja, Valencia, Segunda) was inserted as
# row 46 of the data table, pushing every following record down by one
# row (old row 46 -> new row 47, ..., old row 126 -> new row 127).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 46; Excel shifts rows 46:126 down
# to 47:127 and extends the sheet dimension accordingly.
$ws.Rows.Item(46).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(46, 1).Value = 1
$ws.Cells.Item(46, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(46, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(46, 4).Value = 44973
$ws.Cells.Item(46, 5).Value = 15
$ws.Cells.Item(46, 6).Value = "Fruta"
$ws.Cells.Item(46, 7).Value = 100102
$ws.Cells.Item(46, 8).Value = "Cítricos"
$ws.Cells.Item(46, 9).Value = 100102005
$ws.Cells.Item(46, 10).Value = "Naranja"
$ws.Cells.Item(46, 11).Value = "Valencia"
$ws.Cells.Item(46, 12).Value = "Segunda"
$ws.Cells.Item(46, 13).Value = 300
$ws.Cells.Item(46, 14).Value = 1000
$ws.Cells.Item(46, 15).Value = 1100
$ws.Cells.Item(46, 16).Value = 1050
$ws.Cells.Item(46, 17).Value = "`$/kilo (en caja de 20 kilos)"
$ws.Cells.Item(46, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(46, 19).Value = 1050
$ws.Cells.Item(46, 20).Value = 1
